$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text to bilingual (English/Arabic) labels, keeping each
# cell's original meaning/position.
$ws.Range("A1").Value = "* Date (dd/mm/YYYY)"
$ws.Range("B1").Value = "* Withdrawals/السحوبات"
$ws.Range("C1").Value = "* Deposits/الودائع"
$ws.Range("D1").Value = " * Payee/المستفيد"
$ws.Range("E1").Value = "Description/الوصف"

# Switch the workbook's default font from Calibri to Arial (RTL-friendly).
$wb.Styles.Item("Normal").Font.Name = "Arial"
$wb.Styles.Item("Normal").Name = "عادي"

# Select cell B1 to match the saved view state.
$ws.Range("B1").Select()
